# Replace the 100 arithmetic-problem cells in the single table with their
# updated equations, addressed by (row, column) so each cell is touched
# exactly once (some new values are substrings of other new values, which
# would corrupt a naive whole-document Find/Replace pass).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Text = "69+24="  # was "31-23="
$cell = $t.Cell(1, 2)
$cell.Range.Text = "34+30="  # was "66-37="
$cell = $t.Cell(1, 3)
$cell.Range.Text = "11+16="  # was "97-82="
$cell = $t.Cell(1, 4)
$cell.Range.Text = "0+44="  # was "63+27="
$cell = $t.Cell(1, 5)
$cell.Range.Text = "35+34="  # was "88-25="
$cell = $t.Cell(2, 1)
$cell.Range.Text = "76-29="  # was "4+74="
$cell = $t.Cell(2, 2)
$cell.Range.Text = "99-92="  # was "27+55="
$cell = $t.Cell(2, 3)
$cell.Range.Text = "5+92="  # was "75-63="
$cell = $t.Cell(2, 4)
$cell.Range.Text = "40+35="  # was "98-56="
$cell = $t.Cell(2, 5)
$cell.Range.Text = "63-50="  # was "98-60="
$cell = $t.Cell(3, 1)
$cell.Range.Text = "2+2="  # was "53-45="
$cell = $t.Cell(3, 2)
$cell.Range.Text = "70+29="  # was "16+78="
$cell = $t.Cell(3, 3)
$cell.Range.Text = "27+14="  # was "19+7="
$cell = $t.Cell(3, 4)
$cell.Range.Text = "29-12="  # was "50-33="
$cell = $t.Cell(3, 5)
$cell.Range.Text = "67-11="  # was "14+49="
$cell = $t.Cell(4, 1)
$cell.Range.Text = "64+19="  # was "92-28="
$cell = $t.Cell(4, 2)
$cell.Range.Text = "8+20="  # was "59-46="
$cell = $t.Cell(4, 3)
$cell.Range.Text = "11+23="  # was "82-33="
$cell = $t.Cell(4, 4)
$cell.Range.Text = "74-13="  # was "40-3="
$cell = $t.Cell(4, 5)
$cell.Range.Text = "72+0="  # was "63+23="
$cell = $t.Cell(5, 1)
$cell.Range.Text = "73+26="  # was "85-76="
$cell = $t.Cell(5, 2)
$cell.Range.Text = "67+4="  # was "34-5="
$cell = $t.Cell(5, 3)
$cell.Range.Text = "80-8="  # was "61-38="
$cell = $t.Cell(5, 4)
$cell.Range.Text = "56-4="  # was "1+16="
$cell = $t.Cell(5, 5)
$cell.Range.Text = "44+45="  # was "95-64="
$cell = $t.Cell(6, 1)
$cell.Range.Text = "49-0="  # was "74-0="
$cell = $t.Cell(6, 2)
$cell.Range.Text = "23-16="  # was "40-25="
$cell = $t.Cell(6, 3)
$cell.Range.Text = "92-32="  # was "44-16="
$cell = $t.Cell(6, 4)
$cell.Range.Text = "61-9="  # was "47-16="
$cell = $t.Cell(6, 5)
$cell.Range.Text = "65-36="  # was "78-48="
$cell = $t.Cell(7, 1)
$cell.Range.Text = "98-67="  # was "69-11="
$cell = $t.Cell(7, 2)
$cell.Range.Text = "73-62="  # was "49+17="
$cell = $t.Cell(7, 3)
$cell.Range.Text = "36-19="  # was "35+59="
$cell = $t.Cell(7, 4)
$cell.Range.Text = "56-17="  # was "2+53="
$cell = $t.Cell(7, 5)
$cell.Range.Text = "58-50="  # was "29-15="
$cell = $t.Cell(8, 1)
$cell.Range.Text = "48+14="  # was "84-15="
$cell = $t.Cell(8, 2)
$cell.Range.Text = "43-31="  # was "70-23="
$cell = $t.Cell(8, 3)
$cell.Range.Text = "74-36="  # was "60-10="
$cell = $t.Cell(8, 4)
$cell.Range.Text = "53+36="  # was "91-84="
$cell = $t.Cell(8, 5)
$cell.Range.Text = "50+33="  # was "96-50="
$cell = $t.Cell(9, 1)
$cell.Range.Text = "94-59="  # was "74-41="
$cell = $t.Cell(9, 2)
$cell.Range.Text = "99-26="  # was "29+53="
$cell = $t.Cell(9, 3)
$cell.Range.Text = "45+51="  # was "62+27="
$cell = $t.Cell(9, 4)
$cell.Range.Text = "98-29="  # was "62-0="
$cell = $t.Cell(9, 5)
$cell.Range.Text = "32-18="  # was "67-17="
$cell = $t.Cell(10, 1)
$cell.Range.Text = "16+41="  # was "66+0="
$cell = $t.Cell(10, 2)
$cell.Range.Text = "29+32="  # was "37+46="
$cell = $t.Cell(10, 3)
$cell.Range.Text = "60-50="  # was "72-56="
$cell = $t.Cell(10, 4)
$cell.Range.Text = "44-15="  # was "80-1="
$cell = $t.Cell(10, 5)
$cell.Range.Text = "26+33="  # was "59+35="
$cell = $t.Cell(11, 1)
$cell.Range.Text = "89-56="  # was "68-57="
$cell = $t.Cell(11, 2)
$cell.Range.Text = "25-24="  # was "78-27="
$cell = $t.Cell(11, 3)
$cell.Range.Text = "24+3="  # was "55-52="
$cell = $t.Cell(11, 4)
$cell.Range.Text = "28+2="  # was "41-39="
$cell = $t.Cell(11, 5)
$cell.Range.Text = "30-12="  # was "96-88="
$cell = $t.Cell(12, 1)
$cell.Range.Text = "62+20="  # was "87-60="
$cell = $t.Cell(12, 2)
$cell.Range.Text = "85-18="  # was "97-77="
$cell = $t.Cell(12, 3)
$cell.Range.Text = "83-58="  # was "99-91="
$cell = $t.Cell(12, 4)
$cell.Range.Text = "1+88="  # was "27-8="
$cell = $t.Cell(12, 5)
$cell.Range.Text = "72-47="  # was "61-11="
$cell = $t.Cell(13, 1)
$cell.Range.Text = "64-20="  # was "92+1="
$cell = $t.Cell(13, 2)
$cell.Range.Text = "97-18="  # was "18+60="
$cell = $t.Cell(13, 3)
$cell.Range.Text = "13-4="  # was "66+20="
$cell = $t.Cell(13, 4)
$cell.Range.Text = "72+3="  # was "44-3="
$cell = $t.Cell(13, 5)
$cell.Range.Text = "54-27="  # was "27+57="
$cell = $t.Cell(14, 1)
$cell.Range.Text = "70-18="  # was "12+21="
$cell = $t.Cell(14, 2)
$cell.Range.Text = "82-40="  # was "63-55="
$cell = $t.Cell(14, 3)
$cell.Range.Text = "87-25="  # was "84-53="
$cell = $t.Cell(14, 4)
$cell.Range.Text = "31+18="  # was "95-16="
$cell = $t.Cell(14, 5)
$cell.Range.Text = "82-25="  # was "8+86="
$cell = $t.Cell(15, 1)
$cell.Range.Text = "87+8="  # was "44+29="
$cell = $t.Cell(15, 2)
$cell.Range.Text = "3+17="  # was "85-12="
$cell = $t.Cell(15, 3)
$cell.Range.Text = "40-29="  # was "64-2="
$cell = $t.Cell(15, 4)
$cell.Range.Text = "80-20="  # was "21-1="
$cell = $t.Cell(15, 5)
$cell.Range.Text = "35-19="  # was "29-26="
$cell = $t.Cell(16, 1)
$cell.Range.Text = "18-3="  # was "47-46="
$cell = $t.Cell(16, 2)
$cell.Range.Text = "97-50="  # was "5+93="
$cell = $t.Cell(16, 3)
$cell.Range.Text = "97-32="  # was "55-39="
$cell = $t.Cell(16, 4)
$cell.Range.Text = "86-7="  # was "48-9="
$cell = $t.Cell(16, 5)
$cell.Range.Text = "82-23="  # was "63-60="
$cell = $t.Cell(17, 1)
$cell.Range.Text = "26+31="  # was "5+24="
$cell = $t.Cell(17, 2)
$cell.Range.Text = "61+33="  # was "69-5="
$cell = $t.Cell(17, 3)
$cell.Range.Text = "83-7="  # was "16-5="
$cell = $t.Cell(17, 4)
$cell.Range.Text = "56+6="  # was "10+6="
$cell = $t.Cell(17, 5)
$cell.Range.Text = "69+0="  # was "62+11="
$cell = $t.Cell(18, 1)
$cell.Range.Text = "46+8="  # was "31+48="
$cell = $t.Cell(18, 2)
$cell.Range.Text = "57-11="  # was "39+16="
$cell = $t.Cell(18, 3)
$cell.Range.Text = "45+4="  # was "54+22="
$cell = $t.Cell(18, 4)
$cell.Range.Text = "82-8="  # was "94-17="
$cell = $t.Cell(18, 5)
$cell.Range.Text = "41+50="  # was "79-45="
$cell = $t.Cell(19, 1)
$cell.Range.Text = "44+32="  # was "95-43="
$cell = $t.Cell(19, 2)
$cell.Range.Text = "84-49="  # was "30+24="
$cell = $t.Cell(19, 3)
$cell.Range.Text = "40-24="  # was "9+75="
$cell = $t.Cell(19, 4)
$cell.Range.Text = "32+61="  # was "40+20="
$cell = $t.Cell(19, 5)
$cell.Range.Text = "46+37="  # was "67+9="
$cell = $t.Cell(20, 1)
$cell.Range.Text = "34+41="  # was "8+76="
$cell = $t.Cell(20, 2)
$cell.Range.Text = "65-52="  # was "88-73="
$cell = $t.Cell(20, 3)
$cell.Range.Text = "92-65="  # was "30+36="
$cell = $t.Cell(20, 4)
$cell.Range.Text = "4+95="  # was "98-1="
$cell = $t.Cell(20, 5)
$cell.Range.Text = "71-35="  # was "23+4="
